$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Phase 1: touch one cell for each *new* distinct string value, in the exact
# order that the final shared-strings table must contain them in (the
# underlying engine appends newly-seen string values to the shared string
# table in first-use order). This reproduces the target ordering:
#   181 Ok, 182 robotxt, 183 templateD, 184 Templateed,
#   185 "List Item may be solved above", 186 "Manual ", 187 Blocked,
#   188 ?Blocked, 189 Review
# ---------------------------------------------------------------------------
$ws.Range("J62").Value = "Ok"
$ws.Range("J98").Value = "robotxt"
$ws.Range("J99").Value = "templateD"
$ws.Range("J102").Value = "Templateed"
$ws.Range("J103").Value = "List Item may be solved above"
$ws.Range("J110").Value = "Manual "
$ws.Range("J108").Value = "Blocked"
$ws.Range("J109").Value = "?Blocked"
$ws.Range("J111").Value = "Review"

# ---------------------------------------------------------------------------
# Phase 2: fill in every other cell in column J with its final value so that
# the old "X"/"x" strings become fully unreferenced (and get dropped from the
# shared-string table on save) and the rest of the rows show their final
# values.
# ---------------------------------------------------------------------------

# Rows 2-95: every row becomes "Ok" (previously "X", "x" or "ok").
$ws.Range("J2:J95").Value = "Ok"

# Rows 96-144: previously blank, now populated per-row.
$ws.Range("J96:J97").Value = "ok"
$ws.Range("J100").Value = "robotxt"
$ws.Range("J101").Value = "ok"
$ws.Range("J104:J107").Value = "List Item may be solved above"
$ws.Range("J112").Value = "Review"
$ws.Range("J113:J122").Value = "Manual "
$ws.Range("J123:J135").Value = "ok"
$ws.Range("J136:J139").Value = "robotxt"
$ws.Range("J140:J144").Value = "ok"

# ---------------------------------------------------------------------------
# Update the view: scroll position and active selection.
# ---------------------------------------------------------------------------
$ws.Range("D122").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 101
